$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = [double]"26.731658"
$ws.Cells.Item(2, 8).Value = [double]"80.194974"
$ws.Cells.Item(2, 9).Value = [double]"0.02353393228912"
$ws.Cells.Item(2, 10).Value = [double]"0.02353393228912"
$ws.Cells.Item(2, 13).Value = [double]"7.358839333333333"
$ws.Cells.Item(2, 14).Value = [double]"22.076518"
$ws.Cells.Item(2, 15).Value = [double]"0.3941741583621265"
$ws.Cells.Item(2, 16).Value = [double]"0.3941741583621266"
$ws.Cells.Item(2, 17).Value = [double]"196.7139763356147"
$ws.Cells.Item(2, 18).Value = [double]"1770.425787020532"
$ws.Cells.Item(2, 19).Value = [double]"0.009276467953015148"
$ws.Cells.Item(2, 20).Value = [double]"0.00927646795301515"

$ws.Cells.Item(3, 7).Value = [double]"26.731658"
$ws.Cells.Item(3, 8).Value = [double]"80.194974"
$ws.Cells.Item(3, 9).Value = [double]"0.02353393228912"
$ws.Cells.Item(3, 10).Value = [double]"0.02353393228912"
$ws.Cells.Item(3, 15).Value = [double]"0.1322172945656896"
$ws.Cells.Item(3, 16).Value = [double]"0.1322172945656897"
$ws.Cells.Item(3, 17).Value = [double]"65.98349790972266"
$ws.Cells.Item(3, 18).Value = [double]"593.8514811875041"
$ws.Cells.Item(3, 19).Value = [double]"0.003111592857759573"
$ws.Cells.Item(3, 20).Value = [double]"0.003111592857759574"

$ws.Cells.Item(4, 7).Value = [double]"26.731658"
$ws.Cells.Item(4, 8).Value = [double]"80.194974"
$ws.Cells.Item(4, 9).Value = [double]"0.02353393228912"
$ws.Cells.Item(4, 10).Value = [double]"0.02353393228912"
$ws.Cells.Item(4, 13).Value = [double]"0.4231663333333333"
$ws.Cells.Item(4, 14).Value = [double]"1.269499"
$ws.Cells.Item(4, 15).Value = [double]"0.02266678558034203"
$ws.Cells.Item(4, 16).Value = [double]"0.02266678558034203"
$ws.Cells.Item(4, 17).Value = [double]"11.31193769978067"
$ws.Cells.Item(4, 18).Value = [double]"101.807439298026"
$ws.Cells.Item(4, 19).Value = [double]"0.0005334385970597708"
$ws.Cells.Item(4, 20).Value = [double]"0.0005334385970597709"

$ws.Cells.Item(5, 7).Value = [double]"26.731658"
$ws.Cells.Item(5, 8).Value = [double]"80.194974"
$ws.Cells.Item(5, 9).Value = [double]"0.02353393228912"
$ws.Cells.Item(5, 10).Value = [double]"0.02353393228912"
$ws.Cells.Item(5, 13).Value = [double]"0.133566"
$ws.Cells.Item(5, 14).Value = [double]"0.400698"
$ws.Cells.Item(5, 15).Value = [double]"0.0071544252090564"
$ws.Cells.Item(5, 16).Value = [double]"0.007154425209056401"
$ws.Cells.Item(5, 17).Value = [double]"3.570440632428"
$ws.Cells.Item(5, 18).Value = [double]"32.133965691852"
$ws.Cells.Item(5, 19).Value = [double]"0.0001683717584375065"
$ws.Cells.Item(5, 20).Value = [double]"0.0001683717584375065"

$ws.Cells.Item(6, 7).Value = [double]"26.731658"
$ws.Cells.Item(6, 8).Value = [double]"80.194974"
$ws.Cells.Item(6, 9).Value = [double]"0.02353393228912"
$ws.Cells.Item(6, 10).Value = [double]"0.02353393228912"
$ws.Cells.Item(6, 13).Value = [double]"8.093193666666666"
$ws.Cells.Item(6, 14).Value = [double]"24.279581"
$ws.Cells.Item(6, 15).Value = [double]"0.4335096416047168"
$ws.Cells.Item(6, 16).Value = [double]"0.4335096416047168"
$ws.Cells.Item(6, 17).Value = [double]"216.3444852250993"
$ws.Cells.Item(6, 18).Value = [double]"1947.100367025894"
$ws.Cells.Item(6, 19).Value = [double]"0.01020218655220608"
$ws.Cells.Item(6, 20).Value = [double]"0.01020218655220608"

$ws.Cells.Item(7, 7).Value = [double]"26.731658"
$ws.Cells.Item(7, 8).Value = [double]"80.194974"
$ws.Cells.Item(7, 9).Value = [double]"0.02353393228912"
$ws.Cells.Item(7, 10).Value = [double]"0.02353393228912"
$ws.Cells.Item(7, 13).Value = [double]"0.1918743333333333"
$ws.Cells.Item(7, 14).Value = [double]"0.575623"
$ws.Cells.Item(7, 15).Value = [double]"0.01027769467806845"
$ws.Cells.Item(7, 16).Value = [double]"0.01027769467806845"
$ws.Cells.Item(7, 17).Value = [double]"5.129119057644667"
$ws.Cells.Item(7, 18).Value = [double]"46.162071518802"
$ws.Cells.Item(7, 19).Value = [double]"0.0002418745706419119"
$ws.Cells.Item(7, 20).Value = [double]"0.0002418745706419119"

$ws.Cells.Item(8, 9).Value = [double]"0.9376016087099961"
$ws.Cells.Item(8, 10).Value = [double]"0.9376016087099961"
$ws.Cells.Item(8, 13).Value = [double]"7.358839333333333"
$ws.Cells.Item(8, 14).Value = [double]"22.076518"
$ws.Cells.Item(8, 15).Value = [double]"0.3941741583621265"
$ws.Cells.Item(8, 16).Value = [double]"0.3941741583621266"
$ws.Cells.Item(8, 17).Value = [double]"7837.16628407573"
$ws.Cells.Item(8, 18).Value = [double]"70534.49655668158"
$ws.Cells.Item(8, 19).Value = [double]"0.3695783249922386"
$ws.Cells.Item(8, 20).Value = [double]"0.3695783249922386"

$ws.Cells.Item(9, 9).Value = [double]"0.9376016087099961"
$ws.Cells.Item(9, 10).Value = [double]"0.9376016087099961"
$ws.Cells.Item(9, 15).Value = [double]"0.1322172945656896"
$ws.Cells.Item(9, 16).Value = [double]"0.1322172945656897"
$ws.Cells.Item(9, 19).Value = [double]"0.123967148084074"
$ws.Cells.Item(9, 20).Value = [double]"0.1239671480840741"

$ws.Cells.Item(10, 9).Value = [double]"0.9376016087099961"
$ws.Cells.Item(10, 10).Value = [double]"0.9376016087099961"
$ws.Cells.Item(10, 13).Value = [double]"0.4231663333333333"
$ws.Cells.Item(10, 14).Value = [double]"1.269499"
$ws.Cells.Item(10, 15).Value = [double]"0.02266678558034203"
$ws.Cells.Item(10, 16).Value = [double]"0.02266678558034203"
$ws.Cells.Item(10, 17).Value = [double]"450.6722826701138"
$ws.Cells.Item(10, 18).Value = [double]"4056.050544031024"
$ws.Cells.Item(10, 19).Value = [double]"0.02125241462441323"
$ws.Cells.Item(10, 20).Value = [double]"0.02125241462441323"

$ws.Cells.Item(11, 9).Value = [double]"0.9376016087099961"
$ws.Cells.Item(11, 10).Value = [double]"0.9376016087099961"
$ws.Cells.Item(11, 13).Value = [double]"0.133566"
$ws.Cells.Item(11, 14).Value = [double]"0.400698"
$ws.Cells.Item(11, 15).Value = [double]"0.0071544252090564"
$ws.Cells.Item(11, 16).Value = [double]"0.007154425209056401"
$ws.Cells.Item(11, 17).Value = [double]"142.247833453472"
$ws.Cells.Item(11, 18).Value = [double]"1280.230501081248"
$ws.Cells.Item(11, 19).Value = [double]"0.006708000585406631"
$ws.Cells.Item(11, 20).Value = [double]"0.006708000585406632"

$ws.Cells.Item(12, 9).Value = [double]"0.9376016087099961"
$ws.Cells.Item(12, 10).Value = [double]"0.9376016087099961"
$ws.Cells.Item(12, 13).Value = [double]"8.093193666666666"
$ws.Cells.Item(12, 14).Value = [double]"24.279581"
$ws.Cells.Item(12, 15).Value = [double]"0.4335096416047168"
$ws.Cells.Item(12, 16).Value = [double]"0.4335096416047168"
$ws.Cells.Item(12, 17).Value = [double]"8619.253887985673"
$ws.Cells.Item(12, 18).Value = [double]"77573.28499187106"
$ws.Cells.Item(12, 19).Value = [double]"0.4064593373598763"
$ws.Cells.Item(12, 20).Value = [double]"0.4064593373598764"

$ws.Cells.Item(13, 9).Value = [double]"0.9376016087099961"
$ws.Cells.Item(13, 10).Value = [double]"0.9376016087099961"
$ws.Cells.Item(13, 13).Value = [double]"0.1918743333333333"
$ws.Cells.Item(13, 14).Value = [double]"0.575623"
$ws.Cells.Item(13, 15).Value = [double]"0.01027769467806845"
$ws.Cells.Item(13, 16).Value = [double]"0.01027769467806845"
$ws.Cells.Item(13, 17).Value = [double]"204.3462274231165"
$ws.Cells.Item(13, 18).Value = [double]"1839.116046808048"
$ws.Cells.Item(13, 19).Value = [double]"0.009636383063987147"
$ws.Cells.Item(13, 20).Value = [double]"0.009636383063987147"

$ws.Cells.Item(14, 7).Value = [double]"0.1721486666666666"
$ws.Cells.Item(14, 8).Value = [double]"0.516446"
$ws.Cells.Item(14, 9).Value = [double]"0.0001515556971810586"
$ws.Cells.Item(14, 10).Value = [double]"0.0001515556971810586"
$ws.Cells.Item(14, 13).Value = [double]"7.358839333333333"
$ws.Cells.Item(14, 14).Value = [double]"22.076518"
$ws.Cells.Item(14, 15).Value = [double]"0.3941741583621265"
$ws.Cells.Item(14, 16).Value = [double]"0.3941741583621266"
$ws.Cells.Item(14, 17).Value = [double]"1.266814379447555"
$ws.Cells.Item(14, 18).Value = [double]"11.401329415028"
$ws.Cells.Item(14, 19).Value = [double]"5.973933938132908e-05"
$ws.Cells.Item(14, 20).Value = [double]"5.973933938132908e-05"

$ws.Cells.Item(15, 7).Value = [double]"0.1721486666666666"
$ws.Cells.Item(15, 8).Value = [double]"0.516446"
$ws.Cells.Item(15, 9).Value = [double]"0.0001515556971810586"
$ws.Cells.Item(15, 10).Value = [double]"0.0001515556971810586"
$ws.Cells.Item(15, 15).Value = [double]"0.1322172945656896"
$ws.Cells.Item(15, 16).Value = [double]"0.1322172945656897"
$ws.Cells.Item(15, 17).Value = [double]"0.4249258009795555"
$ws.Cells.Item(15, 18).Value = [double]"3.824332208816"
$ws.Cells.Item(15, 19).Value = [double]"2.003828425729648e-05"
$ws.Cells.Item(15, 20).Value = [double]"2.003828425729648e-05"

$ws.Cells.Item(16, 7).Value = [double]"0.1721486666666666"
$ws.Cells.Item(16, 8).Value = [double]"0.516446"
$ws.Cells.Item(16, 9).Value = [double]"0.0001515556971810586"
$ws.Cells.Item(16, 10).Value = [double]"0.0001515556971810586"
$ws.Cells.Item(16, 13).Value = [double]"0.4231663333333333"
$ws.Cells.Item(16, 14).Value = [double]"1.269499"
$ws.Cells.Item(16, 15).Value = [double]"0.02266678558034203"
$ws.Cells.Item(16, 16).Value = [double]"0.02266678558034203"
$ws.Cells.Item(16, 17).Value = [double]"0.07284752006155554"
$ws.Cells.Item(16, 18).Value = [double]"0.6556276805539999"
$ws.Cells.Item(16, 19).Value = [double]"3.435280491482301e-06"
$ws.Cells.Item(16, 20).Value = [double]"3.435280491482302e-06"

$ws.Cells.Item(17, 7).Value = [double]"0.1721486666666666"
$ws.Cells.Item(17, 8).Value = [double]"0.516446"
$ws.Cells.Item(17, 9).Value = [double]"0.0001515556971810586"
$ws.Cells.Item(17, 10).Value = [double]"0.0001515556971810586"
$ws.Cells.Item(17, 13).Value = [double]"0.133566"
$ws.Cells.Item(17, 14).Value = [double]"0.400698"
$ws.Cells.Item(17, 15).Value = [double]"0.0071544252090564"
$ws.Cells.Item(17, 16).Value = [double]"0.007154425209056401"
$ws.Cells.Item(17, 17).Value = [double]"0.022993208812"
$ws.Cells.Item(17, 18).Value = [double]"0.206938879308"
$ws.Cells.Item(17, 19).Value = [double]"1.084293900488283e-06"
$ws.Cells.Item(17, 20).Value = [double]"1.084293900488283e-06"

$ws.Cells.Item(18, 7).Value = [double]"0.1721486666666666"
$ws.Cells.Item(18, 8).Value = [double]"0.516446"
$ws.Cells.Item(18, 9).Value = [double]"0.0001515556971810586"
$ws.Cells.Item(18, 10).Value = [double]"0.0001515556971810586"
$ws.Cells.Item(18, 13).Value = [double]"8.093193666666666"
$ws.Cells.Item(18, 14).Value = [double]"24.279581"
$ws.Cells.Item(18, 15).Value = [double]"0.4335096416047168"
$ws.Cells.Item(18, 16).Value = [double]"0.4335096416047168"
$ws.Cells.Item(18, 17).Value = [double]"1.393232498791777"
$ws.Cells.Item(18, 18).Value = [double]"12.539092489126"
$ws.Cells.Item(18, 19).Value = [double]"6.570085596811369e-05"
$ws.Cells.Item(18, 20).Value = [double]"6.570085596811369e-05"

$ws.Cells.Item(19, 7).Value = [double]"0.1721486666666666"
$ws.Cells.Item(19, 8).Value = [double]"0.516446"
$ws.Cells.Item(19, 9).Value = [double]"0.0001515556971810586"
$ws.Cells.Item(19, 10).Value = [double]"0.0001515556971810586"
$ws.Cells.Item(19, 13).Value = [double]"0.1918743333333333"
$ws.Cells.Item(19, 14).Value = [double]"0.575623"
$ws.Cells.Item(19, 15).Value = [double]"0.01027769467806845"
$ws.Cells.Item(19, 16).Value = [double]"0.01027769467806845"
$ws.Cells.Item(19, 17).Value = [double]"0.03303091065088889"
$ws.Cells.Item(19, 18).Value = [double]"0.297278195858"
$ws.Cells.Item(19, 19).Value = [double]"1.55764318234872e-06"
$ws.Cells.Item(19, 20).Value = [double]"1.55764318234872e-06"

$ws.Cells.Item(20, 7).Value = [double]"43.33877"
$ws.Cells.Item(20, 8).Value = [double]"130.01631"
$ws.Cells.Item(20, 9).Value = [double]"0.03815444888131313"
$ws.Cells.Item(20, 10).Value = [double]"0.03815444888131313"
$ws.Cells.Item(20, 13).Value = [double]"7.358839333333333"
$ws.Cells.Item(20, 14).Value = [double]"22.076518"
$ws.Cells.Item(20, 15).Value = [double]"0.3941741583621265"
$ws.Cells.Item(20, 16).Value = [double]"0.3941741583621266"
$ws.Cells.Item(20, 17).Value = [double]"318.9230453342867"
$ws.Cells.Item(20, 18).Value = [double]"2870.30740800858"
$ws.Cells.Item(20, 19).Value = [double]"0.01503949777556238"
$ws.Cells.Item(20, 20).Value = [double]"0.01503949777556238"

$ws.Cells.Item(21, 7).Value = [double]"43.33877"
$ws.Cells.Item(21, 8).Value = [double]"130.01631"
$ws.Cells.Item(21, 9).Value = [double]"0.03815444888131313"
$ws.Cells.Item(21, 10).Value = [double]"0.03815444888131313"
$ws.Cells.Item(21, 15).Value = [double]"0.1322172945656896"
$ws.Cells.Item(21, 16).Value = [double]"0.1322172945656897"
$ws.Cells.Item(21, 17).Value = [double]"106.9759174573067"
$ws.Cells.Item(21, 18).Value = [double]"962.7832571157601"
$ws.Cells.Item(21, 19).Value = [double]"0.005044678006732126"
$ws.Cells.Item(21, 20).Value = [double]"0.005044678006732126"

$ws.Cells.Item(22, 7).Value = [double]"43.33877"
$ws.Cells.Item(22, 8).Value = [double]"130.01631"
$ws.Cells.Item(22, 9).Value = [double]"0.03815444888131313"
$ws.Cells.Item(22, 10).Value = [double]"0.03815444888131313"
$ws.Cells.Item(22, 13).Value = [double]"0.4231663333333333"
$ws.Cells.Item(22, 14).Value = [double]"1.269499"
$ws.Cells.Item(22, 15).Value = [double]"0.02266678558034203"
$ws.Cells.Item(22, 16).Value = [double]"0.02266678558034203"
$ws.Cells.Item(22, 17).Value = [double]"18.33950839207667"
$ws.Cells.Item(22, 18).Value = [double]"165.05557552869"
$ws.Cells.Item(22, 19).Value = [double]"0.0008648387117288455"
$ws.Cells.Item(22, 20).Value = [double]"0.0008648387117288456"

$ws.Cells.Item(23, 7).Value = [double]"43.33877"
$ws.Cells.Item(23, 8).Value = [double]"130.01631"
$ws.Cells.Item(23, 9).Value = [double]"0.03815444888131313"
$ws.Cells.Item(23, 10).Value = [double]"0.03815444888131313"
$ws.Cells.Item(23, 13).Value = [double]"0.133566"
$ws.Cells.Item(23, 14).Value = [double]"0.400698"
$ws.Cells.Item(23, 15).Value = [double]"0.0071544252090564"
$ws.Cells.Item(23, 16).Value = [double]"0.007154425209056401"
$ws.Cells.Item(23, 17).Value = [double]"5.78858615382"
$ws.Cells.Item(23, 18).Value = [double]"52.09727538438"
$ws.Cells.Item(23, 19).Value = [double]"0.0002729731509141204"
$ws.Cells.Item(23, 20).Value = [double]"0.0002729731509141204"

$ws.Cells.Item(24, 7).Value = [double]"43.33877"
$ws.Cells.Item(24, 8).Value = [double]"130.01631"
$ws.Cells.Item(24, 9).Value = [double]"0.03815444888131313"
$ws.Cells.Item(24, 10).Value = [double]"0.03815444888131313"
$ws.Cells.Item(24, 13).Value = [double]"8.093193666666666"
$ws.Cells.Item(24, 14).Value = [double]"24.279581"
$ws.Cells.Item(24, 15).Value = [double]"0.4335096416047168"
$ws.Cells.Item(24, 16).Value = [double]"0.4335096416047168"
$ws.Cells.Item(24, 17).Value = [double]"350.7490588851234"
$ws.Cells.Item(24, 18).Value = [double]"3156.74152996611"
$ws.Cells.Item(24, 19).Value = [double]"0.01654032146016354"
$ws.Cells.Item(24, 20).Value = [double]"0.01654032146016354"

$ws.Cells.Item(25, 7).Value = [double]"43.33877"
$ws.Cells.Item(25, 8).Value = [double]"130.01631"
$ws.Cells.Item(25, 9).Value = [double]"0.03815444888131313"
$ws.Cells.Item(25, 10).Value = [double]"0.03815444888131313"
$ws.Cells.Item(25, 13).Value = [double]"0.1918743333333333"
$ws.Cells.Item(25, 14).Value = [double]"0.575623"
$ws.Cells.Item(25, 15).Value = [double]"0.01027769467806845"
$ws.Cells.Item(25, 16).Value = [double]"0.01027769467806845"
$ws.Cells.Item(25, 17).Value = [double]"8.315597601236668"
$ws.Cells.Item(25, 18).Value = [double]"74.84037841113"
$ws.Cells.Item(25, 19).Value = [double]"0.0003921397762121068"
$ws.Cells.Item(25, 20).Value = [double]"0.0003921397762121068"

$ws.Cells.Item(26, 7).Value = [double]"0.4290093333333333"
$ws.Cells.Item(26, 8).Value = [double]"1.287028"
$ws.Cells.Item(26, 9).Value = [double]"0.0003776898762533613"
$ws.Cells.Item(26, 10).Value = [double]"0.0003776898762533613"
$ws.Cells.Item(26, 13).Value = [double]"7.358839333333333"
$ws.Cells.Item(26, 14).Value = [double]"22.076518"
$ws.Cells.Item(26, 15).Value = [double]"0.3941741583621265"
$ws.Cells.Item(26, 16).Value = [double]"0.3941741583621266"
$ws.Cells.Item(26, 17).Value = [double]"3.157010756500444"
$ws.Cells.Item(26, 18).Value = [double]"28.413096808504"
$ws.Cells.Item(26, 19).Value = [double]"0.0001488755890940644"
$ws.Cells.Item(26, 20).Value = [double]"0.0001488755890940644"

$ws.Cells.Item(27, 7).Value = [double]"0.4290093333333333"
$ws.Cells.Item(27, 8).Value = [double]"1.287028"
$ws.Cells.Item(27, 9).Value = [double]"0.0003776898762533613"
$ws.Cells.Item(27, 10).Value = [double]"0.0003776898762533613"
$ws.Cells.Item(27, 15).Value = [double]"0.1322172945656896"
$ws.Cells.Item(27, 16).Value = [double]"0.1322172945656897"
$ws.Cells.Item(27, 17).Value = [double]"1.058951766076444"
$ws.Cells.Item(27, 18).Value = [double]"9.530565894688"
$ws.Cells.Item(27, 19).Value = [double]"4.993713362306955e-05"
$ws.Cells.Item(27, 20).Value = [double]"4.993713362306955e-05"

$ws.Cells.Item(28, 7).Value = [double]"0.4290093333333333"
$ws.Cells.Item(28, 8).Value = [double]"1.287028"
$ws.Cells.Item(28, 9).Value = [double]"0.0003776898762533613"
$ws.Cells.Item(28, 10).Value = [double]"0.0003776898762533613"
$ws.Cells.Item(28, 13).Value = [double]"0.4231663333333333"
$ws.Cells.Item(28, 14).Value = [double]"1.269499"
$ws.Cells.Item(28, 15).Value = [double]"0.02266678558034203"
$ws.Cells.Item(28, 16).Value = [double]"0.02266678558034203"
$ws.Cells.Item(28, 17).Value = [double]"0.1815423065524444"
$ws.Cells.Item(28, 18).Value = [double]"1.633880758972"
$ws.Cells.Item(28, 19).Value = [double]"8.561015440900857e-06"
$ws.Cells.Item(28, 20).Value = [double]"8.561015440900857e-06"

$ws.Cells.Item(29, 7).Value = [double]"0.4290093333333333"
$ws.Cells.Item(29, 8).Value = [double]"1.287028"
$ws.Cells.Item(29, 9).Value = [double]"0.0003776898762533613"
$ws.Cells.Item(29, 10).Value = [double]"0.0003776898762533613"
$ws.Cells.Item(29, 13).Value = [double]"0.133566"
$ws.Cells.Item(29, 14).Value = [double]"0.400698"
$ws.Cells.Item(29, 15).Value = [double]"0.0071544252090564"
$ws.Cells.Item(29, 16).Value = [double]"0.007154425209056401"
$ws.Cells.Item(29, 17).Value = [double]"0.05730106061599999"
$ws.Cells.Item(29, 18).Value = [double]"0.515709545544"
$ws.Cells.Item(29, 19).Value = [double]"2.70215397187244e-06"
$ws.Cells.Item(29, 20).Value = [double]"2.70215397187244e-06"

$ws.Cells.Item(30, 7).Value = [double]"0.4290093333333333"
$ws.Cells.Item(30, 8).Value = [double]"1.287028"
$ws.Cells.Item(30, 9).Value = [double]"0.0003776898762533613"
$ws.Cells.Item(30, 10).Value = [double]"0.0003776898762533613"
$ws.Cells.Item(30, 13).Value = [double]"8.093193666666666"
$ws.Cells.Item(30, 14).Value = [double]"24.279581"
$ws.Cells.Item(30, 15).Value = [double]"0.4335096416047168"
$ws.Cells.Item(30, 16).Value = [double]"0.4335096416047168"
$ws.Cells.Item(30, 17).Value = [double]"3.472055619474222"
$ws.Cells.Item(30, 18).Value = [double]"31.248500575268"
$ws.Cells.Item(30, 19).Value = [double]"0.0001637322028923245"
$ws.Cells.Item(30, 20).Value = [double]"0.0001637322028923245"

$ws.Cells.Item(31, 7).Value = [double]"0.4290093333333333"
$ws.Cells.Item(31, 8).Value = [double]"1.287028"
$ws.Cells.Item(31, 9).Value = [double]"0.0003776898762533613"
$ws.Cells.Item(31, 10).Value = [double]"0.0003776898762533613"
$ws.Cells.Item(31, 13).Value = [double]"0.1918743333333333"
$ws.Cells.Item(31, 14).Value = [double]"0.575623"
$ws.Cells.Item(31, 15).Value = [double]"0.01027769467806845"
$ws.Cells.Item(31, 16).Value = [double]"0.01027769467806845"
$ws.Cells.Item(31, 17).Value = [double]"0.08231587982711111"
$ws.Cells.Item(31, 18).Value = [double]"0.7408429184439999"
$ws.Cells.Item(31, 19).Value = [double]"3.881781231129504e-06"
$ws.Cells.Item(31, 20).Value = [double]"3.881781231129504e-06"

$ws.Cells.Item(32, 7).Value = [double]"0.2053263333333333"
$ws.Cells.Item(32, 8).Value = [double]"0.6159789999999999"
$ws.Cells.Item(32, 9).Value = [double]"0.0001807645461362684"
$ws.Cells.Item(32, 10).Value = [double]"0.0001807645461362684"
$ws.Cells.Item(32, 13).Value = [double]"7.358839333333333"
$ws.Cells.Item(32, 14).Value = [double]"22.076518"
$ws.Cells.Item(32, 15).Value = [double]"0.3941741583621265"
$ws.Cells.Item(32, 16).Value = [double]"0.3941741583621266"
$ws.Cells.Item(32, 17).Value = [double]"1.510963497902444"
$ws.Cells.Item(32, 18).Value = [double]"13.598671481122"
$ws.Cells.Item(32, 19).Value = [double]"7.125271283497539e-05"
$ws.Cells.Item(32, 20).Value = [double]"7.12527128349754e-05"

$ws.Cells.Item(33, 7).Value = [double]"0.2053263333333333"
$ws.Cells.Item(33, 8).Value = [double]"0.6159789999999999"
$ws.Cells.Item(33, 9).Value = [double]"0.0001807645461362684"
$ws.Cells.Item(33, 10).Value = [double]"0.0001807645461362684"
$ws.Cells.Item(33, 15).Value = [double]"0.1322172945656896"
$ws.Cells.Item(33, 16).Value = [double]"0.1322172945656897"
$ws.Cells.Item(33, 17).Value = [double]"0.5068204032204444"
$ws.Cells.Item(33, 18).Value = [double]"4.561383628984"
$ws.Cells.Item(33, 19).Value = [double]"2.39001992435322e-05"
$ws.Cells.Item(33, 20).Value = [double]"2.39001992435322e-05"

$ws.Cells.Item(34, 7).Value = [double]"0.2053263333333333"
$ws.Cells.Item(34, 8).Value = [double]"0.6159789999999999"
$ws.Cells.Item(34, 9).Value = [double]"0.0001807645461362684"
$ws.Cells.Item(34, 10).Value = [double]"0.0001807645461362684"
$ws.Cells.Item(34, 13).Value = [double]"0.4231663333333333"
$ws.Cells.Item(34, 14).Value = [double]"1.269499"
$ws.Cells.Item(34, 15).Value = [double]"0.02266678558034203"
$ws.Cells.Item(34, 16).Value = [double]"0.02266678558034203"
$ws.Cells.Item(34, 17).Value = [double]"0.08688719161344442"
$ws.Cells.Item(34, 18).Value = [double]"0.7819847245209999"
$ws.Cells.Item(34, 19).Value = [double]"4.09735120779864e-06"
$ws.Cells.Item(34, 20).Value = [double]"4.097351207798641e-06"

$ws.Cells.Item(35, 7).Value = [double]"0.2053263333333333"
$ws.Cells.Item(35, 8).Value = [double]"0.6159789999999999"
$ws.Cells.Item(35, 9).Value = [double]"0.0001807645461362684"
$ws.Cells.Item(35, 10).Value = [double]"0.0001807645461362684"
$ws.Cells.Item(35, 13).Value = [double]"0.133566"
$ws.Cells.Item(35, 14).Value = [double]"0.400698"
$ws.Cells.Item(35, 15).Value = [double]"0.0071544252090564"
$ws.Cells.Item(35, 16).Value = [double]"0.007154425209056401"
$ws.Cells.Item(35, 17).Value = [double]"0.02742461703799999"
$ws.Cells.Item(35, 18).Value = [double]"0.246821553342"
$ws.Cells.Item(35, 19).Value = [double]"1.293266425780957e-06"
$ws.Cells.Item(35, 20).Value = [double]"1.293266425780957e-06"

$ws.Cells.Item(36, 7).Value = [double]"0.2053263333333333"
$ws.Cells.Item(36, 8).Value = [double]"0.6159789999999999"
$ws.Cells.Item(36, 9).Value = [double]"0.0001807645461362684"
$ws.Cells.Item(36, 10).Value = [double]"0.0001807645461362684"
$ws.Cells.Item(36, 13).Value = [double]"8.093193666666666"
$ws.Cells.Item(36, 14).Value = [double]"24.279581"
$ws.Cells.Item(36, 15).Value = [double]"0.4335096416047168"
$ws.Cells.Item(36, 16).Value = [double]"0.4335096416047168"
$ws.Cells.Item(36, 17).Value = [double]"1.661745780533222"
$ws.Cells.Item(36, 18).Value = [double]"14.955712024799"
$ws.Cells.Item(36, 19).Value = [double]"7.836317361037301e-05"
$ws.Cells.Item(36, 20).Value = [double]"7.836317361037302e-05"

$ws.Cells.Item(37, 7).Value = [double]"0.2053263333333333"
$ws.Cells.Item(37, 8).Value = [double]"0.6159789999999999"
$ws.Cells.Item(37, 9).Value = [double]"0.0001807645461362684"
$ws.Cells.Item(37, 10).Value = [double]"0.0001807645461362684"
$ws.Cells.Item(37, 13).Value = [double]"0.1918743333333333"
$ws.Cells.Item(37, 14).Value = [double]"0.575623"
$ws.Cells.Item(37, 15).Value = [double]"0.01027769467806845"
$ws.Cells.Item(37, 16).Value = [double]"0.01027769467806845"
$ws.Cells.Item(37, 17).Value = [double]"0.03939685332411111"
$ws.Cells.Item(37, 18).Value = [double]"0.3545716799169999"
$ws.Cells.Item(37, 19).Value = [double]"1.857842813808185e-06"
$ws.Cells.Item(37, 20).Value = [double]"1.857842813808185e-06"

Write-Output "Applied updates to Col1a2-Itgb3 sheet"
